$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (dated 2021-10-20 / serial 44489) was inserted
# above the existing "Especial" row, pushing the old rows 6 and 7 down by
# one row (old row6 -> row7, old row7 -> row8).
$ws.Rows.Item(6).Insert()

# Fill in the newly inserted row 6 with the new record's data.
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C6").Value = "Arica y Parinacota"
$ws.Range("D6").Value = 44489
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100107
$ws.Range("H6").Value = "Otros"
$ws.Range("I6").Value = 100107002
$ws.Range("J6").Value = "Chirimoya"
$ws.Range("K6").Value = "Cultivar IV Región"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 24000
$ws.Range("O6").Value = 25000
$ws.Range("P6").Value = 24500
$ws.Range("Q6").Value = "$/caja 12 kilos"
$ws.Range("R6").Value = "Región de Coquimbo"
$ws.Range("S6").Value = 2042
$ws.Range("T6").Value = 12
